$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 4896.227
$ws.Range("I12").Value = 5345.35
$ws.Range("K12").Value = 5345.35
$ws.Range("M12").Value = -5175.35
$ws.Range("H17").Value = 144542.72
$ws.Range("J17").Value = 144542.72
$ws.Range("L17").Value = 433628.16
$ws.Range("N17").Value = -433964.16
$ws.Range("H70").Value = 6034.1577
$ws.Range("J70").Value = 15950
$ws.Range("L70").Value = 47850
$ws.Range("N70").Value = -48390
$ws.Range("H73").Value = 6034.1577
$ws.Range("J73").Value = 15950
$ws.Range("L73").Value = 47850
$ws.Range("N73").Value = -49722
$ws.Range("H92").Value = 1580
$ws.Range("I92").Value = 1580
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1580
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -332
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 4480.2
$ws.Range("J113").Value = 4903
$ws.Range("L113").Value = 4903
$ws.Range("N113").Value = -11411
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 412.9
$ws.Range("I5").Value = 304.83334
$ws.Range("K5").Value = 304.83334
$ws.Range("M5").Value = -192.83334
$ws.Range("H110").Value = 43133.64
$ws.Range("I110").Value = 46450.086
$ws.Range("K110").Value = 46450.086
$ws.Range("M110").Value = -44405.086
$ws.Range("H132").Value = 3128820.5
$ws.Range("I132").Value = 3128820.5
$ws.Range("K132").Value = 9386461.5
$ws.Range("M132").Value = -9383931.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 412.9
$ws.Range("I4").Value = 304.83334
$ws.Range("K4").Value = 304.83334
$ws.Range("M4").Value = -189.83334
$ws.Range("H105").Value = 2708.25
$ws.Range("I105").Value = 2381
$ws.Range("K105").Value = 2381
$ws.Range("M105").Value = -634
$ws.Range("H107").Value = 167718.17
$ws.Range("I107").Value = 1262
$ws.Range("K107").Value = 1262
$ws.Range("M107").Value = 658
$ws.Range("H134").Value = 83336920
$ws.Range("I134").Value = 125002870
$ws.Range("K134").Value = 375008610
$ws.Range("M134").Value = -375006075

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 29302
$ws.Range("J28").Value = 29302
$ws.Range("L28").Value = 29302
$ws.Range("N28").Value = -29792
$ws.Range("H58").Value = 41677988
$ws.Range("I58").Value = 50013084
$ws.Range("K58").Value = 50013084
$ws.Range("M58").Value = -50012881
$ws.Range("H62").Value = 4672.6
$ws.Range("J62").Value = 4717
$ws.Range("L62").Value = 4717
$ws.Range("N62").Value = -5965
$ws.Range("H65").Value = 4672.6
$ws.Range("J65").Value = 4717
$ws.Range("L65").Value = 23585
$ws.Range("N65").Value = -29825
$ws.Range("H94").Value = 938.6
$ws.Range("I94").Value = 923.625
$ws.Range("K94").Value = 923.625
$ws.Range("M94").Value = -472.625
$ws.Range("H105").Value = 1820883.1
$ws.Range("I105").Value = 2223968.2
$ws.Range("K105").Value = 2223968.2
$ws.Range("M105").Value = -2222221.2
$ws.Range("H136").Value = 41677988
$ws.Range("I136").Value = 50013084
$ws.Range("K136").Value = 150039252
$ws.Range("M136").Value = -150036702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 375
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -131
$ws.Range("H27").Value = 375
$ws.Range("I27").Value = 100
$ws.Range("K27").Value = 300
$ws.Range("M27").Value = -198
$ws.Range("H92").Value = 437
$ws.Range("I92").Value = 449.33334
$ws.Range("K92").Value = 1348.00002
$ws.Range("M92").Value = -100.0000199999999
$ws.Range("H113").Value = 68198.53
$ws.Range("I113").Value = 504999
$ws.Range("J113").Value = 998.46155
$ws.Range("K113").Value = 1514997
$ws.Range("L113").Value = 2995.38465
$ws.Range("M113").Value = -1512827
$ws.Range("N113").Value = -7335.38465
$ws.Range("H132").Value = 1553.0667
$ws.Range("I132").Value = 1731.6666
$ws.Range("J132").Value = 1285.1666
$ws.Range("K132").Value = 15584.9994
$ws.Range("L132").Value = 11566.4994
$ws.Range("M132").Value = -13054.9994
$ws.Range("N132").Value = -16626.4994
$ws.Range("H136").Value = 1094.75
$ws.Range("I136").Value = 1094.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3284.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 1815.75
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I14").Value = 3353389.2
$ws.Range("J14").Value = 60002.5
$ws.Range("K14").Value = 3353389.2
$ws.Range("L14").Value = 60002.5
$ws.Range("M14").Value = -3353221.2
$ws.Range("N14").Value = -60338.5
$ws.Range("H113").Value = 86786.75
$ws.Range("I113").Value = 102344.2
$ws.Range("K113").Value = 102344.2
$ws.Range("M113").Value = -100174.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5471.3335
$ws.Range("I22").Value = 3207.25
$ws.Range("J22").Value = 9999.5
$ws.Range("K22").Value = 3207.25
$ws.Range("L22").Value = 9999.5
$ws.Range("M22").Value = -2912.25
$ws.Range("N22").Value = -10589.5
$ws.Range("H27").Value = 5471.3335
$ws.Range("I27").Value = 3207.25
$ws.Range("J27").Value = 9999.5
$ws.Range("K27").Value = 3207.25
$ws.Range("L27").Value = 9999.5
$ws.Range("M27").Value = -3100.25
$ws.Range("N27").Value = -10213.5
$ws.Range("H40").Value = 1926.5454
$ws.Range("I40").Value = 1829.15
$ws.Range("K40").Value = 1829.15
$ws.Range("M40").Value = -1693.15
$ws.Range("H61").Value = 2931
$ws.Range("I61").Value = 2931
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2931
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2729
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 6251500
$ws.Range("I68").Value = 6251500
$ws.Range("K68").Value = 6251500
$ws.Range("M68").Value = -6250751
$ws.Range("H71").Value = 6251500
$ws.Range("I71").Value = 6251500
$ws.Range("K71").Value = 31257500
$ws.Range("M71").Value = -31253756
$ws.Range("H93").Value = 1399.92
$ws.Range("I93").Value = 1462.8823
$ws.Range("J93").Value = 1266.125
$ws.Range("K93").Value = 1462.8823
$ws.Range("L93").Value = 1266.125
$ws.Range("M93").Value = -214.8823
$ws.Range("N93").Value = -3762.125
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H113").Value = 2931
$ws.Range("I113").Value = 2931
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2931
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -761
$ws.Range("N113").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 10000
$ws.Range("J7").Value = 10000
$ws.Range("L7").Value = 10000
$ws.Range("N7").Value = -10226
$ws.Range("H17").Value = 9995
$ws.Range("J17").Value = 9995
$ws.Range("L17").Value = 9995
$ws.Range("N17").Value = -10339
$ws.Range("H122").Value = 5117.6665
$ws.Range("I122").Value = 5499
$ws.Range("K122").Value = 16497
$ws.Range("M122").Value = -14047
$ws.Range("H132").Value = 15157240
$ws.Range("I132").Value = 26318656
$ws.Range("J132").Value = 9605.857
$ws.Range("K132").Value = 78955968
$ws.Range("L132").Value = 28817.571
$ws.Range("M132").Value = -78953438
$ws.Range("N132").Value = -33877.571

Write-Output "applied 204 cell updates across 8 sheets"